$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain text
# (matching the original inlineStr cell type) -- force text format first
# so the COM layer does not auto-coerce them to numeric cells.
$textForceCells = @('D6', 'D10', 'D15', 'D16', 'D19', 'D21', 'D24', 'D25', 'D27', 'D38', 'D41', 'D45', 'D46', 'D49', 'D50', 'D51')
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = '@'
}

$ws.Range('D2').Value = '27.118.97'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.636.30'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '0.516'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '19.92'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '1.865.34'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').Value = '1.628.96'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').Value = '0.541'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '66.65'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').Value = '27.111.68'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = '216.57'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '6.83'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').Value = '  +3.17%  '
$ws.Range('D24').Value = '9.10'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('D25').Value = '146.52'
$ws.Range('E25').Value = '  -0.30%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = '7.38'
$ws.Range('E27').Value = '  +1.92%  '
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').Value = '1.297.66'
$ws.Range('E34').Value = '  +2.66%  '
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('E36').Value = '  +1.22%  '
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('D38').Value = '0.853'
$ws.Range('E38').Value = '  +2.05%  '
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '0.807'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  +6.25%  '
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('D44').Value = '1.775.55'
$ws.Range('D45').Value = '61.63'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('D46').Value = '91.20'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').Value = '0.0₆0108'
$ws.Range('E48').Value = '  +2.44%  '
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').Value = '7.66'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').Value = '0.0956'
$ws.Range('E51').Value = '  -0.44%  '
